$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'55.841.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.04%  "

# Row 3
$ws.Range("D3").Value = "'2.355.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.71%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "'499.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.44%  "

# Row 6
$ws.Range("D6").Value = "'128.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.80%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "

# Row 8
$ws.Range("D8").Value = "'0.543"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.05%  "

# Row 9
$ws.Range("D9").Value = "'2.357.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.33%  "

# Row 10
$ws.Range("D10").Value = "'0.0977"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.68%  "

# Row 11
$ws.Range("E11").Value = "  -0.07%  "

# Row 12
$ws.Range("D12").Value = "'4.79"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.20%  "

# Row 13
$ws.Range("D13").Value = "'0.322"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "

# Row 14
$ws.Range("D14").Value = "'2.775.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.19%  "

# Row 15
$ws.Range("D15").Value = "'55.846.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.68%  "

# Row 16
$ws.Range("D16").Value = "'21.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.24%  "

# Row 17
$ws.Range("D17").Value = "'0.0000131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.30%  "

# Row 18
$ws.Range("D18").Value = "'2.364.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.88%  "

# Row 19
$ws.Range("D19").Value = "'9.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.81%  "

# Row 20
$ws.Range("D20").Value = "'4.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.85%  "

# Row 21
$ws.Range("D21").Value = "'306.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.02%  "

# Row 22
$ws.Range("D22").Value = "'6.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.81%  "

# Row 23
$ws.Range("E23").Value = "  +0.18%  "

# Row 24
$ws.Range("D24").Value = "'65.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("E26").Value = "  -3.88%  "

# Row 27
$ws.Range("D27").Value = "'0.146"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.61%  "

# Row 28
$ws.Range("D28").Value = "'7.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.01%  "

# Row 29
$ws.Range("D29").Value = "'171.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "

# Row 30
$ws.Range("D30").Value = "'0.0₃0707"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.40%  "

# Row 31
$ws.Range("D31").Value = "'1.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.82%  "

# Row 32
$ws.Range("E32").Value = "  +0.07%  "

# Row 33
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'5.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.89%  "

# Row 35
$ws.Range("D35").Value = "'1.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.92%  "

# Row 36
$ws.Range("D36").Value = "'17.54"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.24%  "

# Row 37
$ws.Range("D37").Value = "'1.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.77%  "

# Row 38
$ws.Range("D38").Value = "'3.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.96%  "

# Row 39
$ws.Range("D39").Value = "'36.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.12%  "

# Row 40
$ws.Range("D40").Value = "'0.789"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.85%  "

# Row 41
$ws.Range("D41").Value = "'1.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.96%  "

# Row 42
$ws.Range("D42").Value = "'3.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.34%  "

# Row 43
$ws.Range("D43").Value = "'127.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.47%  "

# Row 44
$ws.Range("D44").Value = "'4.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.69%  "

# Row 45
$ws.Range("D45").Value = "'0.560"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.32%  "

# Row 46
$ws.Range("D46").Value = "'0.0898"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.55%  "

# Row 47
$ws.Range("D47").Value = "'237.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.25%  "

# Row 48
$ws.Range("D48").Value = "'0.0479"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.33%  "

# Row 49
$ws.Range("D49").Value = "'0.0205"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.88%  "

# Row 50
$ws.Range("D50").Value = "'16.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.11%  "

# Row 51
$ws.Range("E51").Value = "  -1.15%  "
